$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A17").Value = 4
